$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Column C ("Förändrad") bulk update: 45184 -> 45186 for data rows 2..171
$ws.Range("C2:C171").Value2 = 45186

# 2) Add the "Beteckning" text as the second (friendly-name) argument to every
#    single-argument HYPERLINK() formula in columns S,T,U,V,W,X,Y (19..25)
#    for data rows 2..171.
$lastRow = 171
$cols = @(19, 20, 21, 22, 23, 24, 25)

for ($r = 2; $r -le $lastRow; $r++) {
    $beteckning = $ws.Cells.Item($r, 1).Value2
    foreach ($c in $cols) {
        $cell = $ws.Cells.Item($r, $c)
        $f = $cell.Formula
        if ($f -ne "" -and $f.ToUpper().StartsWith("=HYPERLINK(") -and $f.IndexOf(",") -lt 0) {
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
            $cell.Formula = $newFormula
        }
    }
}
